# Word COM-interop edit script
#
# Change 1: In the <语句> ::= ... alternation list, add a new alternative
#           "<continue>|<break>|" right after "...|while循环语句|" and
#           right before "<空>".
# Change 2: Relocate the lone "_GoBack" bookmark so that it sits
#           immediately after the newly inserted text (it previously sat
#           right after the second "<语句>" occurrence, before a line
#           break).  Re-adding a bookmark with the same name simply moves
#           it, so a single Bookmarks.Add call handles both the removal
#           from the old spot and the insertion at the new spot.

$d = $word.ActiveDocument

# --- locate the insertion point -------------------------------------------------
# Unique anchor text: "...while循环语句|" immediately followed by "<空>".
$anchor = $d.Range(0, $d.Content.End)
$found = $anchor.Find.Execute("while循环语句|", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find insertion anchor 'while循环语句|'"
}
$insertAt = $anchor.End

# --- insert the new run's text ---------------------------------------------------
$newText = "<continue>|<break>|"
$insertRange = $d.Range($insertAt, $insertAt)
$insertRange.InsertBefore($newText)

# --- force the freshly-typed text into its own run (same visual formatting
#     as its predecessor run: Times New Roman / sz 15 / szCs 14 / yellow
#     highlight) by toggling a direct-formatting property on and back off.
#     Word (and this host) keep adjoining runs merged while rPr is byte-for
#     -byte identical, so a transient difference forces the split while the
#     final state matches the surrounding run's properties exactly. ---------
$newRange = $d.Range($insertAt, $insertAt + $newText.Length)
$newRange.Bold = 1
$newRange2 = $d.Range($insertAt, $insertAt + $newText.Length)
$newRange2.Bold = 0

# --- move the "_GoBack" bookmark to sit right after the new text ----------------
$bookmarkAt = $insertAt + $newText.Length
$bookmarkRange = $d.Range($bookmarkAt, $bookmarkAt)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
